# Update the OFGEM RHI Domestic figures on Sheet1.
# Source: Data Sources/MANUAL/OFGEM RHIDom.xlsx
# Commit: smart meters, rhi, daily demand, chp

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - East Scotland
$ws.Range("B2").Value = 2112
$ws.Range("C2").Value = 872
$ws.Range("D2").Value = 562
# E2 unchanged (384)

# Row 3 - Highlands & Islands
$ws.Range("B3").Value = 5854
$ws.Range("C3").Value = 1270
$ws.Range("D3").Value = 568
$ws.Range("E3").Value = 540

# Row 4 - North East Scotland
$ws.Range("B4").Value = 2030
# C4, D4, E4 unchanged (341, 305, 175)

# Row 5 - Southern Scotland
# B5 unchanged (4416)
$ws.Range("C5").Value = 1059
$ws.Range("D5").Value = 404
# E5 unchanged (160)

# Row 6 - West Central Scotland
$ws.Range("B6").Value = 770
$ws.Range("C6").Value = 137
# D6, E6 unchanged (73, 30)

# Row 7 totals (F2:F7, B7:F7) recalc automatically via existing formulas.

# Reflect the new selected cell recorded in the workbook view.
$ws.Range("I9").Select()
